$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137, pushing existing rows 137..212 down to 138..213
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new daily price record
$ws.Range("A137").Value = 5
$ws.Range("B137").Value = "Macroferia Regional de Talca"
$ws.Range("C137").Value = "Maule"
$ws.Range("D137").Value = 44452
$ws.Range("E137").Value = 7
$ws.Range("F137").Value = 100112043
$ws.Range("G137").Value = "Pepino ensalada"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 300
$ws.Range("K137").Value = 15000
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = 15000
$ws.Range("N137").Value = "`$/caja 60 unidades"
$ws.Range("O137").Value = "Región de Arica y Parinacota"
$ws.Range("P137").Value = 250
$ws.Range("Q137").Value = 60
$ws.Range("R137").Value = "Hortaliza"

# D column holds dates formatted as date/time; make sure the new cell matches
# the number format used by the rest of the column (style carried over from
# the Insert() already, but set explicitly for safety).
$ws.Range("D137").NumberFormat = $ws.Range("D138").NumberFormat
